$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text-like values (safe to assign directly; Excel will not
# reinterpret them as numbers because of formatting/extra dots/spaces).
$plainValues = @{
    "D2" = "67.875.88"
    "E2" = "  +3.80%  "
    "D3" = "3.313.00"
    "E3" = "  +1.52%  "
    "E4" = "  -0.43%  "
    "E5" = "  +0.44%  "
    "E6" = "  -0.50%  "
    "E7" = "  +0.11%  "
    "E8" = "  +3.80%  "
    "D9" = "3.308.56"
    "E9" = "  +1.36%  "
    "E10" = "  +2.32%  "
    "E11" = "  +2.20%  "
    "E12" = "  +0.75%  "
    "E13" = "  +4.54%  "
    "E14" = "  +14.27%  "
    "D15" = "3.843.19"
    "E15" = "  +1.55%  "
    "E16" = "  +1.41%  "
    "D17" = "67.940.95"
    "E17" = "  +3.82%  "
    "E18" = "  +1.56%  "
    "D19" = "3.310.80"
    "E19" = "  +1.23%  "
    "E20" = "  +0.21%  "
    "E21" = "  +1.01%  "
    "E22" = "  +2.41%  "
    "E23" = "  -6.54%  "
    "E24" = "  +5.80%  "
    "E25" = "  +1.44%  "
    "E26" = "  +0.93%  "
    "E27" = "  +2.63%  "
    "E28" = "  +1.47%  "
    "E29" = "  +10.54%  "
    "E30" = "  +2.87%  "
    "E31" = "  +6.30%  "
    "E32" = "  +7.09%  "
    "D33" = "3.924.06"
    "E33" = "  +4.75%  "
    "E34" = "  +1.72%  "
    "E35" = "  +2.27%  "
    "E36" = "  +0.24%  "
    "E37" = "  -5.82%  "
    "E38" = "  -0.24%  "
    "E39" = "  +4.49%  "
    "E40" = "  +2.33%  "
    "E41" = "  +4.01%  "
    "B42" = "InjectiveProtocol"
    "C42" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "E42" = "  +0.91%  "
    "B43" = "PEPE"
    "C43" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D43" = "0.0₃0687"
    "E43" = "  +3.53%  "
    "E44" = "  -1.15%  "
    "E45" = "  +2.22%  "
    "E46" = "  +3.22%  "
    "E47" = "  +3.20%  "
    "E48" = "  +11.21%  "
    "E49" = "  +0.23%  "
    "E50" = "  +3.69%  "
    "E51" = "  +0.69%  "
}
foreach ($cellref in $plainValues.Keys) {
    $ws.Range($cellref).Value = $plainValues[$cellref]
}

# Numeric-looking text values (e.g. "1.00", "0.104") must be forced to
# stay text so trailing zeros / exact formatting survive, matching the
# original inline-string cells (which carried no numeric style).
$textValues = @{
    "D4" = "0.997"
    "D5" = "579.08"
    "D6" = "180.11"
    "D7" = "1.00"
    "D12" = "46.03"
    "D14" = "701.23"
    "D16" = "8.41"
    "D21" = "10.87"
    "D23" = "16.95"
    "D24" = "5.18"
    "D25" = "99.13"
    "D27" = "2.75"
    "D28" = "9.40"
    "D29" = "33.45"
    "D30" = "8.51"
    "D31" = "6.76"
    "D32" = "583.04"
    "D34" = "10.91"
    "D35" = "0.104"
    "D38" = "55.42"
    "D39" = "0.131"
    "D40" = "3.18"
    "D41" = "2.64"
    "D42" = "32.42"
    "D44" = "3.35"
    "D46" = "0.0413"
    "D51" = "128.64"
}
foreach ($cellref in $textValues.Keys) {
    $cell = $ws.Range($cellref)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$cellref]
    $cell.Style = "Normal"
}
